$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.26
$ws.Range("A3").Value = -21.575
$ws.Range("A14").Value = -21.606
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.536
$ws.Range("A25").Value = -21.521
$ws.Range("B25").Value = 6.843999999999999
$ws.Range("A26").Value = -21.351
$ws.Range("B27").Value = 6.529999999999999
$ws.Range("A29").Value = -21.219
$ws.Range("B31").Value = 6.415999999999999
$ws.Range("B39").Value = 7.582000000000001
$ws.Range("B48").Value = 5.274
$ws.Range("B51").Value = 6.226
$ws.Range("B52").Value = 5.813000000000001
$ws.Range("A53").Value = -22.01
$ws.Range("B55").Value = 4.670999999999999
$ws.Range("B56").Value = 5.003
$ws.Range("A57").Value = -21.607
$ws.Range("B57").Value = 6.210000000000001
$ws.Range("A59").Value = -22.5
$ws.Range("A69").Value = -21.507
$ws.Range("B73").Value = 6.804
$ws.Range("A79").Value = -21.246
$ws.Range("A83").Value = -22.006
$ws.Range("B89").Value = 5.787999999999999
$ws.Range("B90").Value = 5.767
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 5.899
$ws.Range("A93").Value = -21.439
